{"js": "// Remove the trailing \"Requisitos\" section (its Heading2 heading paragraph\n// plus the following ListBullet paragraph naming the prerequisite course)\n// that was appended at the very end of the document body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet headingIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].style === \"Heading 2\" && items[i].text.trim() === \"Requisitos\") {\n    headingIndex = i;\n    break;\n  }\n}\n\nif (headingIndex !== -1) {\n  // Delete the paragraph right after the heading first (if present) so the\n  // heading's own index/range stays valid, then delete the heading itself.\n  if (headingIndex + 1 < items.length) {\n    items[headingIndex + 1].delete();\n  }\n  items[headingIndex].delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Find the \"Requisitos\" Heading2 paragraph, then remove it together with the\n# paragraph right after it (the bulleted prerequisite course entry) -- this\n# deletes the whole \"Requisitos\" section that was appended at the end of the\n# course description.\n$headingIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Style.NameLocal -eq \"Heading 2\" -and $p.Range.Text.Trim() -eq \"Requisitos\") {\n        $headingIndex = $i\n        break\n    }\n}\n\nif ($headingIndex -ne -1) {\n    $endIndex = $headingIndex\n    if ($d.Paragraphs.Count -gt $headingIndex) {\n        $endIndex = $headingIndex + 1\n    }\n    $startPara = $d.Paragraphs.Item($headingIndex)\n    $endPara = $d.Paragraphs.Item($endIndex)\n    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)\n    $rng.Delete()\n}\n"}
